# Add a new "Yearly demand" worksheet at the end of the workbook, matching
# the layout/formatting already used by the other hourly-profile sheets
# (header row 0-23 across B:Y, row labels 0/1/2 down A2:A4), then fill in
# the yearly-demand figures.

$wb = $excel.ActiveWorkbook

# Clone the formatting/layout of an existing hourly-profile sheet so the new
# sheet picks up the same header style (bold, bordered, centered) without
# having to hand-build the style table.
$template = $wb.Worksheets.Item("DG Dispatch")
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Yearly demand"

# Row 2 (A2 = 0)
$ws.Range("B2").Value = -32.5
$ws.Range("C2").Value = -19.5
$ws.Range("D2").Value = -13
$ws.Range("E2").Value = -13
$ws.Range("F2").Value = -13
$ws.Range("G2").Value = 142.5
$ws.Range("H2").Value = 291.5
$ws.Range("I2").Value = 327
$ws.Range("J2").Value = 388.5
$ws.Range("K2").Value = 502
$ws.Range("L2").Value = 596
$ws.Range("M2").Value = 670.5
$ws.Range("N2").Value = 745
$ws.Range("O2").Value = 651
$ws.Range("P2").Value = 576.5
$ws.Range("Q2").Value = 502
$ws.Range("R2").Value = 320.5
$ws.Range("S2").Value = 139
$ws.Range("T2").Value = 32
$ws.Range("U2").Value = -117
$ws.Range("V2").Value = -97.5
$ws.Range("W2").Value = -78
$ws.Range("X2").Value = -52
$ws.Range("Y2").Value = -39

# Row 3 (A3 = 1)
$ws.Range("B3").Value = -32.5
$ws.Range("C3").Value = -19.5
$ws.Range("D3").Value = -13
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = -19.5
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 324
$ws.Range("J3").Value = 486
$ws.Range("K3").Value = 648
$ws.Range("L3").Value = 729
$ws.Range("M3").Value = 751.5
$ws.Range("N3").Value = 583
$ws.Range("O3").Value = 567
$ws.Range("P3").Value = 333.5
$ws.Range("Q3").Value = 340
$ws.Range("R3").Value = 243
$ws.Range("S3").Value = 57.99999999999999
$ws.Range("T3").Value = -130
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = -78
$ws.Range("X3").Value = 0
$ws.Range("Y3").Value = -39

# Row 4 (A4 = 2)
$ws.Range("B4").Value = -32.5
$ws.Range("C4").Value = -19.5
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = -19.5
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 81
$ws.Range("K4").Value = 324
$ws.Range("L4").Value = 567
$ws.Range("M4").Value = 589.5
$ws.Range("N4").Value = 648
$ws.Range("O4").Value = 567
$ws.Range("P4").Value = 324
$ws.Range("Q4").Value = 162
$ws.Range("R4").Value = 81
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = -130
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0
$ws.Range("X4").Value = 0
$ws.Range("Y4").Value = -39

# Copying a sheet makes the copy the active tab in real Excel too; restore
# the original active sheet so the view state is left untouched.
$template.Select()
